# Flops tracker update: refresh percentage-used figures (column H) and
# remove the scratch "runs for testing" rows (60-61) now that the check
# is complete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated percentage-used values from the latest experiment report.
$ws.Range("H32").Value = 0.96509999999999996
$ws.Range("H36").Value = 0.80989999999999995
$ws.Range("H39").Value = 0.73509999999999998
$ws.Range("H40").Value = 0.65149999999999997
$ws.Range("H42").Value = 0.6744
$ws.Range("H43").Value = 0.61619999999999997
$ws.Range("H45").Value = 0.72289999999999999

# Remove the "amount of runs I can have for testing" scratch rows.
$ws.Range("A60").ClearContents()
$ws.Range("B60").ClearContents()
$ws.Range("B61").ClearContents()

# Move the active selection to reflect where the author was last working.
$ws.Range("H46").Select() | Out-Null
